# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The account-statement (Estado de Cuenta) data block is refreshed:
#   - the totals (VALOR MORA / Cant. Trabajadores / Cant. Periodos) are
#     reset to reflect a single remaining worker,
#   - the first data row is rewritten with that worker's info and updated
#     amounts,
#   - the other worker rows are removed entirely,
#   - the Nombre Trabajador column is narrowed to fit the shorter content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# VALOR MORA total
$ws.Range("E11").Value = 40000

# Cant. Trabajadores / Cant. Periodos
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Remaining worker row (ALEXANDER OLIVERA ROCHA) with refreshed amounts
$ws.Range("C16").Value = "1143357977"
$ws.Range("D16").Value = "ALEXANDER OLIVERA ROCHA"
$ws.Range("E16").Value = "2209"
$ws.Range("F16").Value = 40000
$ws.Range("G16").Value = 2000000

# Drop the other worker rows (their data has been superseded)
$ws.Rows("17:20").Delete()

# Column D no longer needs to fit the longer names that were removed
$ws.Columns("D:D").ColumnWidth = 26.1666666667
